# Auto-generated script applying market-data refresh updates
# (scheduled runner price/profit recalculation) to Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 779
$ws.Range("I12").Value = 686.8
$ws.Range("K12").Value = 686.8
$ws.Range("M12").Value = -516.8
# Row 33
$ws.Range("H33").Value = 206.86667
$ws.Range("I33").Value = 227.33333
$ws.Range("K33").Value = 227.33333
$ws.Range("M33").Value = 1.666670000000011
# Row 100
$ws.Range("H100").Value = 1798.4445
$ws.Range("J100").Value = 2472.75
$ws.Range("L100").Value = 2472.75
$ws.Range("N100").Value = -3554.75
# Row 140
$ws.Range("H140").Value = 48264.535
$ws.Range("J140").Value = 48264.535
$ws.Range("L140").Value = 48264.535
$ws.Range("N140").Value = -58624.535

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2840.0557
$ws.Range("I32").Value = 2077.0698
$ws.Range("K32").Value = 2077.0698
$ws.Range("M32").Value = -1790.0698
# Row 45
$ws.Range("H45").Value = 2064.818
$ws.Range("I45").Value = 926
$ws.Range("J45").Value = 2715.5715
$ws.Range("K45").Value = 926
$ws.Range("L45").Value = 2715.5715
$ws.Range("M45").Value = -549
$ws.Range("N45").Value = -3469.5715
# Row 61
$ws.Range("H61").Value = 4908.5713
$ws.Range("I61").Value = 3016.5
$ws.Range("J61").Value = 7431.3335
$ws.Range("K61").Value = 3016.5
$ws.Range("L61").Value = 7431.3335
$ws.Range("M61").Value = -2804.5
$ws.Range("N61").Value = -7855.3335
# Row 88
$ws.Range("H88").Value = 2746.2144
$ws.Range("I88").Value = 2049.6667
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 2049.6667
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = -1643.6667
$ws.Range("N88").Value = -4812
# Row 91
$ws.Range("H91").Value = 2746.2144
$ws.Range("I91").Value = 2049.6667
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 2049.6667
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = -645.6667000000002
$ws.Range("N91").Value = -6808
# Row 97
$ws.Range("H97").Value = 486
$ws.Range("I97").Value = 492.1
$ws.Range("J97").Value = 465.66666
$ws.Range("K97").Value = 492.1
$ws.Range("L97").Value = 465.66666
$ws.Range("M97").Value = 3.899999999999977
$ws.Range("N97").Value = -1457.66666
# Row 122
$ws.Range("H122").Value = 1971.75
$ws.Range("I122").Value = 2061.3635
$ws.Range("K122").Value = 6184.0905
$ws.Range("M122").Value = -3734.0905
# Row 136
$ws.Range("H136").Value = 4908.5713
$ws.Range("I136").Value = 3016.5
$ws.Range("J136").Value = 7431.3335
$ws.Range("K136").Value = 9049.5
$ws.Range("L136").Value = 22294.0005
$ws.Range("M136").Value = -6499.5
$ws.Range("N136").Value = -27394.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 917.0833
$ws.Range("I94").Value = 729.5
$ws.Range("J94").Value = 1855
$ws.Range("K94").Value = 729.5
$ws.Range("L94").Value = 1855
$ws.Range("M94").Value = -278.5
$ws.Range("N94").Value = -2757
# Row 105
$ws.Range("H105").Value = 2485.5
$ws.Range("I105").Value = 2492.077
$ws.Range("K105").Value = 2492.077
$ws.Range("M105").Value = -745.0770000000002

$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 574.2
$ws.Range("I105").Value = 574.2
$ws.Range("K105").Value = 574.2
$ws.Range("M105").Value = 1172.8
# Row 132
$ws.Range("H132").Value = 3598.0833
$ws.Range("I132").Value = 1593
$ws.Range("K132").Value = 4779
$ws.Range("M132").Value = -2249

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 599.2308
$ws.Range("I5").Value = 519.1
$ws.Range("K5").Value = 1557.3
$ws.Range("M5").Value = -1445.3
# Row 50
$ws.Range("H50").Value = 116980.664
$ws.Range("I50").Value = 232791.33
$ws.Range("J50").Value = 1170
$ws.Range("K50").Value = 698373.99
$ws.Range("L50").Value = 3510
$ws.Range("M50").Value = -697892.99
$ws.Range("N50").Value = -4472
# Row 53
$ws.Range("H53").Value = 116980.664
$ws.Range("I53").Value = 232791.33
$ws.Range("J53").Value = 1170
$ws.Range("K53").Value = 698373.99
$ws.Range("L53").Value = 3510
$ws.Range("M53").Value = -697892.99
$ws.Range("N53").Value = -4472
# Row 61
$ws.Range("H61").Value = 285
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 285
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 855
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1285
# Row 131
$ws.Range("H131").Value = 8078112
$ws.Range("J131").Value = 15266.109
$ws.Range("L131").Value = 45798.327
$ws.Range("N131").Value = -55878.327
# Row 133
$ws.Range("H133").Value = 16670175
$ws.Range("I133").Value = 62500956
$ws.Range("K133").Value = 187502868
$ws.Range("M133").Value = -187497808
# Row 135
$ws.Range("H135").Value = 599.2308
$ws.Range("I135").Value = 519.1
$ws.Range("K135").Value = 4671.900000000001
$ws.Range("M135").Value = -2136.900000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
# Row 132
$ws.Range("H132").Value = 3777.158
$ws.Range("I132").Value = 2971.7273
$ws.Range("J132").Value = 4884.625
$ws.Range("K132").Value = 8915.1819
$ws.Range("L132").Value = 14653.875
$ws.Range("M132").Value = -6385.1819
$ws.Range("N132").Value = -19713.875

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1994
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 1994
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 1994
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -4490
# Row 122
$ws.Range("H122").Value = 5672.087
$ws.Range("I122").Value = 2605.889
$ws.Range("J122").Value = 7643.2144
$ws.Range("K122").Value = 7817.667
$ws.Range("L122").Value = 22929.6432
$ws.Range("M122").Value = -5367.667
$ws.Range("N122").Value = -27829.6432

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30630
# Row 73
$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -32184
# Row 122
$ws.Range("H122").Value = 21784.56
$ws.Range("I122").Value = 32951.375
$ws.Range("K122").Value = 98854.125
$ws.Range("M122").Value = -96404.125

Write-Output "Applied market data refresh."